$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Preço Atual" (column C) values for each commodity row
$ws.Range("C2").Value = 68.68000000000001
$ws.Range("C3").Value = 133.17
$ws.Range("C4").Value = 279.15
$ws.Range("C5").Value = 392.46
$ws.Range("C6").Value = 410.01
$ws.Range("C7").Value = 147.82
$ws.Range("C8").Value = 1092.76
$ws.Range("C9").Value = 324.22
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 9.609999999999999

# Update "Comprar" (column D) boolean values that changed
$ws.Range("D7").Value = $false
$ws.Range("D8").Value = $true
$ws.Range("D11").Value = $false
